# backup before deadline cleaning
# Update a handful of recomputed result values in row 2 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 2.7210000000000001
$ws.Range("K2").Value = 0.71399999999999997
$ws.Range("Y2").Value = 1.4980000000000002
$ws.Range("AA2").Value = 2.1080000000000005
$ws.Range("AE2").Value = 1.2809999999999999
$ws.Range("AF2").Value = 1.159
$ws.Range("AV2").Value = 2.8940000000000001
$ws.Range("BN2").Value = 2.5410000000000004
$ws.Range("BY2").Value = 1.3689999999999998
$ws.Range("CQ2").Value = 1.2969999999999999
